$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns remain plain text (matching the original inline-string cells)
# so numeric-looking values (trailing zeros, etc.) are preserved exactly.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "37.038.62"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").Value = "2.022.72"
$ws.Range("E3").Value = "  -2.83%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "227.21"
$ws.Range("E5").Value = "  -2.59%  "
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  -4.02%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "55.25"
$ws.Range("E8").Value = "  -4.55%  "
$ws.Range("D9").Value = "0.381"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("D10").Value = "0.0796"
$ws.Range("E10").Value = "  +2.18%  "
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  -3.69%  "
$ws.Range("D12").Value = "2.321.93"
$ws.Range("E12").Value = "  -2.89%  "
$ws.Range("D13").Value = "14.34"
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("D14").Value = "20.69"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").Value = "0.747"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").Value = "5.17"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("D17").Value = "2.024.36"
$ws.Range("E17").Value = "  -3.22%  "
$ws.Range("D18").Value = "36.962.12"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "6.13"
$ws.Range("E19").Value = "  +1.55%  "
$ws.Range("D20").Value = "68.94"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "0.0₃0831"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "226.69"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").Value = "2.28"
$ws.Range("E25").Value = "  -4.27%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.35"
$ws.Range("E26").Value = "  -3.32%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "167.32"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("E28").Value = "  -4.09%  "
$ws.Range("D29").Value = "18.80"
$ws.Range("E29").Value = "  -3.59%  "
$ws.Range("D30").Value = "1.34"
$ws.Range("E30").Value = "  -2.58%  "
$ws.Range("E31").Value = "  -3.96%  "
$ws.Range("D32").Value = "4.50"
$ws.Range("E32").Value = "  -3.18%  "
$ws.Range("D33").Value = "0.0613"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("D34").Value = "4.45"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("D35").Value = "2.38"
$ws.Range("E35").Value = "  -3.97%  "
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "3.19"
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D39").Value = "5.43"
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("D40").Value = "0.0220"
$ws.Range("E40").Value = "  -5.10%  "
$ws.Range("D41").Value = "1.499.42"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("D42").Value = "17.02"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("D43").Value = "0.0931"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("D44").Value = "95.35"
$ws.Range("E44").Value = "  -5.15%  "
$ws.Range("E45").Value = "  -2.65%  "
$ws.Range("E46").Value = "  -5.05%  "
$ws.Range("D47").Value = "7.28"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  -4.10%  "
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("D50").Value = "3.72"
$ws.Range("E50").Value = "  -6.76%  "
$ws.Range("D51").Value = "2.209.82"
$ws.Range("E51").Value = "  -2.85%  "
